$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd15c50a8e0>),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=0.001,
                                                                 class_weight=''balanced'',
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver=''liblinear''),
                                    n_estimators=10, random_state=42))])'
$ws.Range("B2").Value = 0.7189743589743589
$ws.Range("C2").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd15c4787c0>, ''scaler'': None, ''model__n_estimators'': 10, ''model__estimator__solver'': ''liblinear'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': ''balanced'', ''model__estimator__C'': 0.001}'
$ws.Range("D2").Value = 0.461189340887457
$ws.Range("E2").Value = 0.4208105018974584
$ws.Range("F2").Value = 0.823529411764706
$ws.Range("G2").Value = 0.4611512305818457
$ws.Range("H2").Value = 0.418478835978836
$ws.Range("I2").Value = 0.7777777777777778
$ws.Range("J2").Value = 0.4820844896700586
$ws.Range("K2").Value = 0.4480676328502415
$ws.Range("L2").Value = 0.875
$ws.Range("M2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range("N2").Value = '[1 1 1 1 1 1 1 0 0 1 0 1 1 1 0 0 0 1 1 1 1 1 1 1]'
$ws.Range("O2").Value = 42

$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd15c70c1f0>),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=3,
                                                                 class_weight=''balanced'',
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver=''liblinear''),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B3").Value = 0.6796037296037296
$ws.Range("C3").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd15c1913a0>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__solver'': ''liblinear'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': ''balanced'', ''model__estimator__C'': 3}'
$ws.Range("D3").Value = 0.4584153659534136
$ws.Range("E3").Value = 0.4032062301627519
$ws.Range("F3").Value = 0.7567567567567567
$ws.Range("G3").Value = 0.4451953854377775
$ws.Range("H3").Value = 0.4289803312629399
$ws.Range("I3").Value = 0.6666666666666666
$ws.Range("J3").Value = 0.5004316990440949
$ws.Range("K3").Value = 0.4275362318840579
$ws.Range("L3").Value = 0.875
$ws.Range("M3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range("N3").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 0 1 1 1 1 1 0 1 0 1 1]'
$ws.Range("O3").Value = 69

$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd15c1915e0>),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=3,
                                                                 class_weight=''balanced'',
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver=''saga''),
                                    random_state=42))])'
$ws.Range("B4").Value = 0.6660317460317461
$ws.Range("C4").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd15c713730>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 50, ''model__estimator__solver'': ''saga'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': ''balanced'', ''model__estimator__C'': 3}'
$ws.Range("D4").Value = 0.4249604970844431
$ws.Range("E4").Value = 0.3715995896464646
$ws.Range("F4").Value = 0.6206896551724138
$ws.Range("G4").Value = 0.4277908729255168
$ws.Range("H4").Value = 0.3822414434523809
$ws.Range("I4").Value = 0.9
$ws.Range("J4").Value = 0.4386805555555555
$ws.Range("K4").Value = 0.3975
$ws.Range("L4").Value = 0.4736842105263158
$ws.Range("M4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("N4").Value = '[0 1 1 0 1 1 0 1 1 0 0 0 0 1 1 0 0 1 0 0 0 0 1 0]'
$ws.Range("O4").Value = 23

$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd21f672760>),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=0.0001,
                                                                 class_weight=''balanced'',
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver=''liblinear''),
                                    n_estimators=10, random_state=42))])'
$ws.Range("B5").Value = 0.7599999999999999
$ws.Range("C5").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd15c0bae50>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 10, ''model__estimator__solver'': ''liblinear'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': ''balanced'', ''model__estimator__C'': 0.0001}'
$ws.Range("D5").Value = 0.4184476714198709
$ws.Range("E5").Value = 0.3749191882798439
$ws.Range("F5").Value = 0.7368421052631579
$ws.Range("G5").Value = 0.4405589646078986
$ws.Range("H5").Value = 0.383127764767109
$ws.Range("I5").Value = 0.5833333333333334
$ws.Range("J5").Value = 0.414519906323185
$ws.Range("K5").Value = 0.3912568306010928
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]'
$ws.Range("N5").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("O5").Value = 99

$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fd15c14a730>),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=0.0001,
                                                                 class_weight=''balanced'',
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver=''liblinear''),
                                    random_state=42))])'
$ws.Range("B6").Value = 0.7081573981573981
$ws.Range("C6").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fd15c214490>, ''scaler'': StandardScaler(), ''model__n_estimators'': 50, ''model__estimator__solver'': ''liblinear'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': ''balanced'', ''model__estimator__C'': 0.0001}'
$ws.Range("D6").Value = 0.4800822517027923
$ws.Range("E6").Value = 0.3994435415330937
$ws.Range("F6").Value = 0.6428571428571428
$ws.Range("G6").Value = 0.5011317084361564
$ws.Range("H6").Value = 0.4133848614072494
$ws.Range("I6").Value = 0.5294117647058824
$ws.Range("J6").Value = 0.4837256027554535
$ws.Range("K6").Value = 0.418407960199005
$ws.Range("L6").Value = 0.8181818181818182
$ws.Range("M6").Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range("N6").Value = '[1 1 1 0 1 1 1 0 1 0 0 1 0 1 1 0 1 1 0 1 1 1 1 1]'
$ws.Range("O6").Value = 89
